$d = $word.ActiveDocument

# Locate the "Should I include other features to predict price?" paragraph
# (last item of the "Questions:" list). The two new question paragraphs are
# inserted right after it, before the "Findings:" heading.
$findRange = $d.Content
$found = $findRange.Find.Execute("Should I include other features to predict price?",
                                  $true, $false, $false, $false, $false,
                                  $true, 1, $false, "", 0)

$anchorPara = $findRange.Paragraphs(1)
$rng = $anchorPara.Range

# Insert the first new paragraph. InsertParagraphAfter() clones the anchor
# paragraph's formatting (Heading1 style, numId 7 list, spacing, run props),
# producing an empty paragraph right after it; then fill in its text.
$rng.InsertParagraphAfter()
$p1 = $d.Range($rng.End, $rng.End)
$p1.Text = "Does dropout improve performance?"

# Re-locate the paragraph we just created to get a fresh (non-stale) Range
# whose End sits right at its paragraph mark, then insert the second new
# paragraph after it the same way.
$p1para = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Does dropout improve performance?*") {
        $p1para = $p
    }
}
$rng2 = $p1para.Range
$rng2.InsertParagraphAfter()
$p2 = $d.Range($rng2.End, $rng2.End)
$p2.Text = "Does multiple layers improve performance?"
